$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "69.302.64"
Set-TextValue $ws.Range("E2") "  -0.02%  "

Set-TextValue $ws.Range("D3") "3.673.30"
Set-TextValue $ws.Range("E3") "  -0.46%  "

Set-TextValue $ws.Range("E4") "  +0.09%  "

Set-TextValue $ws.Range("D5") "682.91"
Set-TextValue $ws.Range("E5") "  -0.02%  "

Set-TextValue $ws.Range("D6") "157.59"
Set-TextValue $ws.Range("E6") "  -3.16%  "

Set-TextValue $ws.Range("E7") "  +0.06%  "

Set-TextValue $ws.Range("D8") "0.493"
Set-TextValue $ws.Range("E8") "  -1.17%  "

Set-TextValue $ws.Range("D9") "0.145"
Set-TextValue $ws.Range("E9") "  -2.22%  "

Set-TextValue $ws.Range("D10") "6.98"
Set-TextValue $ws.Range("E10") "  -4.02%  "

Set-TextValue $ws.Range("D11") "0.434"
Set-TextValue $ws.Range("E11") "  -3.80%  "

Set-TextValue $ws.Range("D12") "0.0000231"
Set-TextValue $ws.Range("E12") "  -2.79%  "

Set-TextValue $ws.Range("D13") "4.297.29"
Set-TextValue $ws.Range("E13") "  -0.30%  "

Set-TextValue $ws.Range("D14") "32.06"
Set-TextValue $ws.Range("E14") "  -4.77%  "

Set-TextValue $ws.Range("D15") "3.687.98"
Set-TextValue $ws.Range("E15") "  -0.04%  "

Set-TextValue $ws.Range("D16") "69.359.14"
Set-TextValue $ws.Range("E16") "  -0.03%  "

Set-TextValue $ws.Range("E17") "  +2.00%  "

Set-TextValue $ws.Range("D18") "15.74"
Set-TextValue $ws.Range("E18") "  -3.70%  "

Set-TextValue $ws.Range("D19") "6.35"
Set-TextValue $ws.Range("E19") "  -4.73%  "

Set-TextValue $ws.Range("D20") "472.16"
Set-TextValue $ws.Range("E20") "  -2.18%  "

Set-TextValue $ws.Range("D21") "9.91"
Set-TextValue $ws.Range("E21") "  +1.21%  "

Set-TextValue $ws.Range("D22") "0.645"
Set-TextValue $ws.Range("E22") "  -3.42%  "

Set-TextValue $ws.Range("D23") "79.91"

Set-TextValue $ws.Range("D24") "3.823.74"
Set-TextValue $ws.Range("E24") "  -0.29%  "

Set-TextValue $ws.Range("E25") "  -0.08%  "

Set-TextValue $ws.Range("D26") "0.0000121"
Set-TextValue $ws.Range("E26") "  -5.91%  "

Set-TextValue $ws.Range("D27") "10.87"
Set-TextValue $ws.Range("E27") "  -5.87%  "

Set-TextValue $ws.Range("D28") "9.08"
Set-TextValue $ws.Range("E28") "  -5.62%  "

Set-TextValue $ws.Range("D29") "2.69"
Set-TextValue $ws.Range("E29") "  -2.19%  "

Set-TextValue $ws.Range("D30") "1.73"
Set-TextValue $ws.Range("E30") "  -5.40%  "

Set-TextValue $ws.Range("E31") "  +0.03%  "

Set-TextValue $ws.Range("D32") "6.50"
Set-TextValue $ws.Range("E32") "  -5.03%  "

Set-TextValue $ws.Range("B33") "EthereumClassic"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D33") "26.79"
Set-TextValue $ws.Range("E33") "  -1.18%  "

Set-TextValue $ws.Range("B34") "ImmutableX"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D34") "1.97"
Set-TextValue $ws.Range("E34") "  -7.27%  "

Set-TextValue $ws.Range("D35") "3.655.56"
Set-TextValue $ws.Range("E35") "  +0.10%  "

Set-TextValue $ws.Range("D36") "0.157"
Set-TextValue $ws.Range("E36") "  -5.26%  "

Set-TextValue $ws.Range("D37") "8.13"
Set-TextValue $ws.Range("E37") "  -5.07%  "

Set-TextValue $ws.Range("D38") "6.04"
Set-TextValue $ws.Range("E38") "  -1.09%  "

Set-TextValue $ws.Range("D40") "2.19"
Set-TextValue $ws.Range("E40") "  +0.13%  "

Set-TextValue $ws.Range("D41") "0.0895"
Set-TextValue $ws.Range("E41") "  -5.22%  "

Set-TextValue $ws.Range("E42") "  -0.02%  "

Set-TextValue $ws.Range("D43") "0.938"
Set-TextValue $ws.Range("E43") "  -2.20%  "

Set-TextValue $ws.Range("D44") "165.69"
Set-TextValue $ws.Range("E44") "  +4.91%  "

Set-TextValue $ws.Range("D45") "47.54"
Set-TextValue $ws.Range("E45") "  -1.30%  "

Set-TextValue $ws.Range("D46") "0.000278"
Set-TextValue $ws.Range("E46") "  -1.31%  "

Set-TextValue $ws.Range("D47") "2.68"
Set-TextValue $ws.Range("E47") "  -6.15%  "

Set-TextValue $ws.Range("D48") "1.10"
Set-TextValue $ws.Range("E48") "  +1.73%  "

Set-TextValue $ws.Range("D49") "1.27"
Set-TextValue $ws.Range("E49") "  -2.79%  "

Set-TextValue $ws.Range("D50") "7.73"
Set-TextValue $ws.Range("E50") "  -4.63%  "

Set-TextValue $ws.Range("D51") "26.78"
Set-TextValue $ws.Range("E51") "  -3.88%  "
